# Generate Report for Handoff
#
# Adds a new file "f3ee664a-882b-4597-b10b-13db7f9b9078.md" (status
# "Ready for handoff") to the localization status report, and flips the
# existing "e9c99ed0-7199-40a0-9c39-deb40e4f5ce6.md" entry to its
# successor file "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md" (now also
# "Ready for handoff" instead of "In Translation"), across all three
# worksheets (Overview, zh-cn, de-de). The ".localization-config" row
# shifts down by one row on every sheet to make room.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276   # BGR-int for RGB(0x64,0x95,0xED) - the workbook's HyperLink font color
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Set-HyperlinkLook($range) {
    $range.Font.Color = $hyperlinkColor
    $range.Font.Underline = $true
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# This engine's Hyperlinks.Delete() on any range clears every hyperlink
# on the sheet, and there is no reliable in-place hyperlink edit, so we
# drop them all and recreate the full, correct set after the cell values
# are in place.
$ws.Range("A1").Hyperlinks.Delete()

# Push ".localization-config" (old row 4) down to row 5.
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# Row 3: e9c99ed0-... -> 80ed4e75-..., now "Ready for handoff".
$ws.Range("A3").Value = "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

# Row 4: new entry for f3ee664a-..., also "Ready for handoff".
$ws.Range("A4").Value = "f3ee664a-882b-4597-b10b-13db7f9b9078.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

Set-HyperlinkLook $ws.Range("A3")
Set-HyperlinkLook $ws.Range("A4")
Set-HyperlinkLook $ws.Range("A5")

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/e2e/ad6d125a-fc5e-4c40-a197-417a1450ae7f.md", "", "", "ad6d125a-fc5e-4c40-a197-417a1450ae7f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/e2e/80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md", "", "", "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/e2e/f3ee664a-882b-4597-b10b-13db7f9b9078.md", "", "", "f3ee664a-882b-4597-b10b-13db7f9b9078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A1").Hyperlinks.Delete()

# Push ".localization-config" (old row 4) down to row 5.
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("D5").NumberFormat = $dateFormat
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Row 3: e9c99ed0-... -> 80ed4e75-..., now "Ready for handoff".
$ws.Range("A3").Value = "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.96556f1f0ced9539490ebebbe60e1fb52eab807b.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-03 07:10:45"
$ws.Range("D3").NumberFormat = $dateFormat
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

# Row 4: new entry for f3ee664a-...
$ws.Range("A4").Value = "f3ee664a-882b-4597-b10b-13db7f9b9078.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "f3ee664a-882b-4597-b10b-13db7f9b9078.bc1621d6e13c85a400c171bc17f122f14ff728f0.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-03 07:10:45"
$ws.Range("D4").NumberFormat = $dateFormat
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

Set-HyperlinkLook $ws.Range("A3")
Set-HyperlinkLook $ws.Range("C3")
Set-HyperlinkLook $ws.Range("A4")
Set-HyperlinkLook $ws.Range("C4")
Set-HyperlinkLook $ws.Range("A5")

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/e2e/ad6d125a-fc5e-4c40-a197-417a1450ae7f.md", "", "", "ad6d125a-fc5e-4c40-a197-417a1450ae7f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85803a5924198649af61450fabb05973708b5158/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ad6d125a-fc5e-4c40-a197-417a1450ae7f.2eb5086cfa2a94e4014402a4731cc59511f277da.zh-cn.xlf", "", "", "ad6d125a-fc5e-4c40-a197-417a1450ae7f.2eb5086cfa2a94e4014402a4731cc59511f277da.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/e2e/80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md", "", "", "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85803a5924198649af61450fabb05973708b5158/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.96556f1f0ced9539490ebebbe60e1fb52eab807b.zh-cn.xlf", "", "", "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.96556f1f0ced9539490ebebbe60e1fb52eab807b.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/e2e/f3ee664a-882b-4597-b10b-13db7f9b9078.md", "", "", "f3ee664a-882b-4597-b10b-13db7f9b9078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/85803a5924198649af61450fabb05973708b5158/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f3ee664a-882b-4597-b10b-13db7f9b9078.bc1621d6e13c85a400c171bc17f122f14ff728f0.zh-cn.xlf", "", "", "f3ee664a-882b-4597-b10b-13db7f9b9078.bc1621d6e13c85a400c171bc17f122f14ff728f0.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A1").Hyperlinks.Delete()

# Push ".localization-config" (old row 4) down to row 5.
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("D5").NumberFormat = $dateFormat
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# Row 3: e9c99ed0-... -> 80ed4e75-..., now "Ready for handoff".
$ws.Range("A3").Value = "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.96556f1f0ced9539490ebebbe60e1fb52eab807b.de-de.xlf"
$ws.Range("D3").Value = "2016-03-03 07:10:55"
$ws.Range("D3").NumberFormat = $dateFormat
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

# Row 4: new entry for f3ee664a-...
$ws.Range("A4").Value = "f3ee664a-882b-4597-b10b-13db7f9b9078.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "f3ee664a-882b-4597-b10b-13db7f9b9078.bc1621d6e13c85a400c171bc17f122f14ff728f0.de-de.xlf"
$ws.Range("D4").Value = "2016-03-03 07:10:55"
$ws.Range("D4").NumberFormat = $dateFormat
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

Set-HyperlinkLook $ws.Range("A3")
Set-HyperlinkLook $ws.Range("C3")
Set-HyperlinkLook $ws.Range("A4")
Set-HyperlinkLook $ws.Range("C4")
Set-HyperlinkLook $ws.Range("A5")

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/e2e/ad6d125a-fc5e-4c40-a197-417a1450ae7f.md", "", "", "ad6d125a-fc5e-4c40-a197-417a1450ae7f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3ad16d61b0ee185ee92292fa29cf98136d92511/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ad6d125a-fc5e-4c40-a197-417a1450ae7f.2eb5086cfa2a94e4014402a4731cc59511f277da.de-de.xlf", "", "", "ad6d125a-fc5e-4c40-a197-417a1450ae7f.2eb5086cfa2a94e4014402a4731cc59511f277da.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/e2e/80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md", "", "", "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3ad16d61b0ee185ee92292fa29cf98136d92511/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.96556f1f0ced9539490ebebbe60e1fb52eab807b.de-de.xlf", "", "", "80ed4e75-bfe1-40ab-ab92-394fd7b6dd08.96556f1f0ced9539490ebebbe60e1fb52eab807b.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/e2e/f3ee664a-882b-4597-b10b-13db7f9b9078.md", "", "", "f3ee664a-882b-4597-b10b-13db7f9b9078.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3ad16d61b0ee185ee92292fa29cf98136d92511/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f3ee664a-882b-4597-b10b-13db7f9b9078.bc1621d6e13c85a400c171bc17f122f14ff728f0.de-de.xlf", "", "", "f3ee664a-882b-4597-b10b-13db7f9b9078.bc1621d6e13c85a400c171bc17f122f14ff728f0.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/cc36dd5bced1bb7b4587fd9ac48a7787427cb1ef/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "done"
